$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45; existing rows 45..95 shift down to 46..96
$ws.Rows(45).Insert()

# Populate the newly inserted row 45 with the new record's data
$ws.Range("A45").Value = 9
$ws.Range("B45").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C45").Value = "Metropolitana"
$ws.Range("D45").Value = 44669
$ws.Range("E45").Value = 13
$ws.Range("F45").Value = 100114007
$ws.Range("G45").Value = "Jengibre"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 610
$ws.Range("K45").Value = 10000
$ws.Range("L45").Value = 11000
$ws.Range("M45").Value = 10500
$ws.Range("N45").Value = "$/caja 13 kilos"
$ws.Range("O45").Value = "Perú"
$ws.Range("P45").Value = 808
$ws.Range("Q45").Value = 13
$ws.Range("R45").Value = "Hortaliza"
